$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Header et footer" table cell — merge the two runs ("Header et "
# + "footer") that were split apart by a spell-check proofErr marker back
# into a single run with the plain text "Header et footer".
# ---------------------------------------------------------------------------
$table = $d.Tables.Item(1)
$foundCell = $null
for ($ri = 1; $ri -le $table.Rows.Count; $ri++) {
    for ($ci = 1; $ci -le $table.Columns.Count; $ci++) {
        try {
            $cell = $table.Cell($ri, $ci)
        } catch {
            continue
        }
        if ($cell.Range.Text -like "Header et footer*") {
            $foundCell = $cell
        }
    }
}
if ($foundCell -ne $null) {
    $cellRange = $d.Range($foundCell.Range.Start, $foundCell.Range.End)
    $cellRange.Text = "Header et footer"
}

# ---------------------------------------------------------------------------
# Change 2: final paragraph of the main body ("Fonctionnement du site"
# section) currently just holds the placeholder text "p". Flesh it out with
# the real paragraph about the header/footer behaviour, add a first-line
# indent, then append a blank paragraph and a second paragraph (tab-indented)
# describing the footer.
# ---------------------------------------------------------------------------
$bodyEnd = $d.Content.End
$placeholder = $d.Range($bodyEnd - 2, $bodyEnd - 1)

$text1 = "Lorsque vous serez sur la page d’accueil, vous aurez en haut de votre écran un header avec deux menus, un afficher et un en burger. Si vous cliquez sur l’un des noms afficher cela vous amène sur la page. Si vous cliquez sur le burger, cela vous ouvrira un menu a partir duquel vous pouvez naviguer entre les pages comme avec le menu afficher sur la page d’accueil."
$text2 = "Lorsque vous faites défilé la page d’accueil vous avez toute la page et à la fin vous avez a nouveau un menu dans le footer."

# First-line indent (540 twips = 27 pt) on the paragraph.
$placeholder.ParagraphFormat.FirstLineIndent = 27

# Replace the "p" placeholder text with the real paragraph text.
$placeholder.Text = $text1

# Append a blank paragraph, then the tab-indented footer paragraph, right
# after the paragraph we just filled in.
$newEnd = $d.Content.End
$tail = $d.Range($newEnd - 1, $newEnd - 1)
$tail.InsertAfter("" + [char]13 + [char]13 + [char]9 + $text2)

Write-Output "done"
